$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.413.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.792.31'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.42'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5475'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3826'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07569'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.45'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.120'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.09'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.184'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.386'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +6.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.795.71'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06463'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.952'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.385.45'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.10%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.121'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.90'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.71'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.403'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.001.02'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.85'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.121'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1026'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.727'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.699'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2306'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +14.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06521'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +9.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02318'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.180'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.769'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +9.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.60'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6383'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.45%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.158'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.387'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.52'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5956'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.674'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.989'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.146'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06897'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.94%  '

Write-Output "Updated 91 cells"
